$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "MuSCs" rows (originally rows 8 and 9).
# Deleting the same row index twice removes both, shifting later rows up.
$ws.Rows(8).EntireRow.Delete()
$ws.Rows(8).EntireRow.Delete()

# Refresh the recomputed TPM-derived statistics for every remaining data row.
# Row 2: ECs -> ECs
$ws.Cells.Item(2,7).Value = 0.3176433333333333
$ws.Cells.Item(2,8).Value = 0.95293
$ws.Cells.Item(2,9).Value = 0.3649408890199488
$ws.Cells.Item(2,10).Value = 0.3649408890199488
$ws.Cells.Item(2,13).Value = 0.2689956666666667
$ws.Cells.Item(2,14).Value = 0.806987
$ws.Cells.Item(2,15).Value = 0.1265890356442505
$ws.Cells.Item(2,16).Value = 0.1265890356442505
$ws.Cells.Item(2,17).Value = 0.08544468021222223
$ws.Cells.Item(2,18).Value = 0.7690021219100001
$ws.Cells.Item(2,19).Value = 0.04619751520819078
$ws.Cells.Item(2,20).Value = 0.04619751520819078

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,7).Value = 0.3176433333333333
$ws.Cells.Item(3,8).Value = 0.95293
$ws.Cells.Item(3,9).Value = 0.3649408890199488
$ws.Cells.Item(3,10).Value = 0.3649408890199488
$ws.Cells.Item(3,13).Value = 1.855956666666666
$ws.Cells.Item(3,14).Value = 5.567869999999999
$ws.Cells.Item(3,15).Value = 0.8734109643557494
$ws.Cells.Item(3,16).Value = 0.8734109643557494
$ws.Cells.Item(3,17).Value = 0.5895322621222222
$ws.Cells.Item(3,18).Value = 5.3057903591
$ws.Cells.Item(3,19).Value = 0.318743373811758
$ws.Cells.Item(3,20).Value = 0.318743373811758

# Row 4: FAPs -> ECs
$ws.Cells.Item(4,7).Value = 0.07444366666666667
$ws.Cells.Item(4,8).Value = 0.223331
$ws.Cells.Item(4,9).Value = 0.08552843722593914
$ws.Cells.Item(4,10).Value = 0.08552843722593914
$ws.Cells.Item(4,13).Value = 0.2689956666666667
$ws.Cells.Item(4,14).Value = 0.806987
$ws.Cells.Item(4,15).Value = 0.1265890356442505
$ws.Cells.Item(4,16).Value = 0.1265890356442505
$ws.Cells.Item(4,17).Value = 0.02002502374411111
$ws.Cells.Item(4,18).Value = 0.180225213697
$ws.Cells.Item(4,19).Value = 0.01082696238859145
$ws.Cells.Item(4,20).Value = 0.01082696238859145

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5,7).Value = 0.07444366666666667
$ws.Cells.Item(5,8).Value = 0.223331
$ws.Cells.Item(5,9).Value = 0.08552843722593914
$ws.Cells.Item(5,10).Value = 0.08552843722593914
$ws.Cells.Item(5,13).Value = 1.855956666666666
$ws.Cells.Item(5,14).Value = 5.567869999999999
$ws.Cells.Item(5,15).Value = 0.8734109643557494
$ws.Cells.Item(5,16).Value = 0.8734109643557494
$ws.Cells.Item(5,17).Value = 0.1381642194411111
$ws.Cells.Item(5,18).Value = 1.24347797497
$ws.Cells.Item(5,19).Value = 0.07470147483734768
$ws.Cells.Item(5,20).Value = 0.07470147483734768

# Row 6: Inflammatory-Mac -> ECs
$ws.Cells.Item(6,7).Value = 0.311492
$ws.Cells.Item(6,8).Value = 0.934476
$ws.Cells.Item(6,9).Value = 0.3578736131801976
$ws.Cells.Item(6,10).Value = 0.3578736131801976
$ws.Cells.Item(6,13).Value = 0.2689956666666667
$ws.Cells.Item(6,14).Value = 0.806987
$ws.Cells.Item(6,15).Value = 0.1265890356442505
$ws.Cells.Item(6,16).Value = 0.1265890356442505
$ws.Cells.Item(6,17).Value = 0.08378999820133334
$ws.Cells.Item(6,18).Value = 0.754109983812
$ws.Cells.Item(6,19).Value = 0.04530287557500476
$ws.Cells.Item(6,20).Value = 0.04530287557500476

# Row 7: Inflammatory-Mac -> FAPs
$ws.Cells.Item(7,7).Value = 0.311492
$ws.Cells.Item(7,8).Value = 0.934476
$ws.Cells.Item(7,9).Value = 0.3578736131801976
$ws.Cells.Item(7,10).Value = 0.3578736131801976
$ws.Cells.Item(7,13).Value = 1.855956666666666
$ws.Cells.Item(7,14).Value = 5.567869999999999
$ws.Cells.Item(7,15).Value = 0.8734109643557494
$ws.Cells.Item(7,16).Value = 0.8734109643557494
$ws.Cells.Item(7,17).Value = 0.5781156540133332
$ws.Cells.Item(7,18).Value = 5.203040886119999
$ws.Cells.Item(7,19).Value = 0.3125707376051928
$ws.Cells.Item(7,20).Value = 0.3125707376051928

# Row 8: Neutrophils -> ECs
$ws.Cells.Item(8,7).Value = 0.08275166666666667
$ws.Cells.Item(8,8).Value = 0.248255
$ws.Cells.Item(8,9).Value = 0.09507351054500056
$ws.Cells.Item(8,10).Value = 0.09507351054500057
$ws.Cells.Item(8,13).Value = 0.2689956666666667
$ws.Cells.Item(8,14).Value = 0.806987
$ws.Cells.Item(8,15).Value = 0.1265890356442505
$ws.Cells.Item(8,16).Value = 0.1265890356442505
$ws.Cells.Item(8,17).Value = 0.02225983974277778
$ws.Cells.Item(8,18).Value = 0.200338557685
$ws.Cells.Item(8,19).Value = 0.0120352640152051
$ws.Cells.Item(8,20).Value = 0.01203526401520511

# Row 9: Neutrophils -> FAPs
$ws.Cells.Item(9,7).Value = 0.08275166666666667
$ws.Cells.Item(9,8).Value = 0.248255
$ws.Cells.Item(9,9).Value = 0.09507351054500056
$ws.Cells.Item(9,10).Value = 0.09507351054500057
$ws.Cells.Item(9,13).Value = 1.855956666666666
$ws.Cells.Item(9,14).Value = 5.567869999999999
$ws.Cells.Item(9,15).Value = 0.8734109643557494
$ws.Cells.Item(9,16).Value = 0.8734109643557494
$ws.Cells.Item(9,17).Value = 0.1535835074277778
$ws.Cells.Item(9,18).Value = 1.38225156685
$ws.Cells.Item(9,19).Value = 0.08303824652979544
$ws.Cells.Item(9,20).Value = 0.08303824652979545

# Row 10: Resolving-Mac -> ECs
$ws.Cells.Item(10,7).Value = 0.08406599999999999
$ws.Cells.Item(10,8).Value = 0.252198
$ws.Cells.Item(10,9).Value = 0.09658355002891401
$ws.Cells.Item(10,10).Value = 0.09658355002891403
$ws.Cells.Item(10,13).Value = 0.2689956666666667
$ws.Cells.Item(10,14).Value = 0.806987
$ws.Cells.Item(10,15).Value = 0.1265890356442505
$ws.Cells.Item(10,16).Value = 0.1265890356442505
$ws.Cells.Item(10,17).Value = 0.022613389714
$ws.Cells.Item(10,18).Value = 0.203520507426
$ws.Cells.Item(10,19).Value = 0.01222641845725845
$ws.Cells.Item(10,20).Value = 0.01222641845725845

# Row 11: Resolving-Mac -> FAPs
$ws.Cells.Item(11,7).Value = 0.08406599999999999
$ws.Cells.Item(11,8).Value = 0.252198
$ws.Cells.Item(11,9).Value = 0.09658355002891401
$ws.Cells.Item(11,10).Value = 0.09658355002891403
$ws.Cells.Item(11,13).Value = 1.855956666666666
$ws.Cells.Item(11,14).Value = 5.567869999999999
$ws.Cells.Item(11,15).Value = 0.8734109643557494
$ws.Cells.Item(11,16).Value = 0.8734109643557494
$ws.Cells.Item(11,17).Value = 0.15602285314
$ws.Cells.Item(11,18).Value = 1.40420567826
$ws.Cells.Item(11,19).Value = 0.08435713157165556
$ws.Cells.Item(11,20).Value = 0.08435713157165557

